$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 149.75
$ws.Range("I12").Value = 149.75
$ws.Range("K12").Value = 149.75
$ws.Range("M12").Value = 20.25
$ws.Range("H32").Value = 1200.3334
$ws.Range("J32").Value = 2002
$ws.Range("L32").Value = 2002
$ws.Range("N32").Value = -2654
$ws.Range("H40").Value = 8327.727999999999
$ws.Range("I40").Value = 10576.25
$ws.Range("K40").Value = 10576.25
$ws.Range("M40").Value = -10401.25
$ws.Range("H43").Value = 10608.2
$ws.Range("I43").Value = 13090.363
$ws.Range("K43").Value = 13090.363
$ws.Range("M43").Value = -13021.363

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 17699.5
$ws.Range("I22").Value = 1331.6666
$ws.Range("K22").Value = 1331.6666
$ws.Range("M22").Value = -1032.6666
$ws.Range("H45").Value = 5688.609
$ws.Range("I45").Value = 6269.55
$ws.Range("J45").Value = 1815.6666
$ws.Range("K45").Value = 6269.55
$ws.Range("L45").Value = 1815.6666
$ws.Range("M45").Value = -5892.55
$ws.Range("N45").Value = -2569.6666
$ws.Range("H53").Value = 45000
$ws.Range("J53").Value = 45000
$ws.Range("L53").Value = 45000
$ws.Range("N53").Value = -46364
$ws.Range("H92").Value = 33000
$ws.Range("J92").Value = 33000
$ws.Range("L92").Value = 33000
$ws.Range("N92").Value = -37992
$ws.Range("H97").Value = 1988.909
$ws.Range("I97").Value = 1355.579
$ws.Range("K97").Value = 1355.579
$ws.Range("M97").Value = -859.579
$ws.Range("H102").Value = 4149.7144
$ws.Range("I102").Value = 3812.5
$ws.Range("K102").Value = 3812.5
$ws.Range("M102").Value = -2190.5
$ws.Range("H122").Value = 1417.909
$ws.Range("I122").Value = 1329.4445
$ws.Range("K122").Value = 3988.3335
$ws.Range("M122").Value = -1538.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 15003.333
$ws.Range("J19").Value = 17505
$ws.Range("L19").Value = 17505
$ws.Range("N19").Value = -17851
$ws.Range("H20").Value = 4951.1724
$ws.Range("I20").Value = 4418.3687
$ws.Range("K20").Value = 4418.3687
$ws.Range("M20").Value = -4171.3687
$ws.Range("H50").Value = 20555
$ws.Range("J50").Value = 20555
$ws.Range("L50").Value = 20555
$ws.Range("N50").Value = -21703
$ws.Range("H80").Value = 385
$ws.Range("J80").Value = 343.14285
$ws.Range("L80").Value = 343.14285
$ws.Range("N80").Value = -2339.14285
$ws.Range("H83").Value = 385
$ws.Range("J83").Value = 343.14285
$ws.Range("L83").Value = 1715.71425
$ws.Range("N83").Value = -11699.71425
$ws.Range("H86").Value = 4370.654
$ws.Range("I86").Value = 4519.9414
$ws.Range("K86").Value = 4519.9414
$ws.Range("M86").Value = -3396.9414
$ws.Range("H89").Value = 4370.654
$ws.Range("I89").Value = 4519.9414
$ws.Range("K89").Value = 22599.707
$ws.Range("M89").Value = -16983.707
$ws.Range("H102").Value = 33137
$ws.Range("I102").Value = 14899.5
$ws.Range("K102").Value = 14899.5
$ws.Range("M102").Value = -11654.5
$ws.Range("H105").Value = 2079.8928
$ws.Range("I105").Value = 2150.5
$ws.Range("J105").Value = 1903.375
$ws.Range("K105").Value = 2150.5
$ws.Range("L105").Value = 1903.375
$ws.Range("M105").Value = -403.5
$ws.Range("N105").Value = -5397.375
$ws.Range("H134").Value = 2240.2307
$ws.Range("I134").Value = 2202.0518
$ws.Range("J134").Value = 2556.5715
$ws.Range("K134").Value = 6606.155400000001
$ws.Range("L134").Value = 7669.7145
$ws.Range("M134").Value = -4071.155400000001
$ws.Range("N134").Value = -12739.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4924.909
$ws.Range("I58").Value = 4691
$ws.Range("K58").Value = 4691
$ws.Range("M58").Value = -4488
$ws.Range("H99").Value = 2711.3
$ws.Range("I99").Value = 2157.5715
$ws.Range("K99").Value = 2157.5715
$ws.Range("M99").Value = -659.5715
$ws.Range("H105").Value = 1214.8
$ws.Range("I105").Value = 1306.8572
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1306.8572
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 440.1428000000001
$ws.Range("N105").Value = -4494
$ws.Range("H107").Value = 881.14813
$ws.Range("I107").Value = 631.4091
$ws.Range("K107").Value = 631.4091
$ws.Range("M107").Value = 1288.5909
$ws.Range("H122").Value = 1602.56
$ws.Range("I122").Value = 1271.1333
$ws.Range("J122").Value = 2099.7
$ws.Range("K122").Value = 3813.3999
$ws.Range("L122").Value = 6299.099999999999
$ws.Range("M122").Value = -1363.3999
$ws.Range("N122").Value = -11199.1
$ws.Range("H126").Value = 2711.3
$ws.Range("I126").Value = 2157.5715
$ws.Range("K126").Value = 6472.7145
$ws.Range("M126").Value = -4002.7145
$ws.Range("H135").Value = 49536.25
$ws.Range("J135").Value = 49536.25
$ws.Range("L135").Value = 49536.25
$ws.Range("N135").Value = -59676.25
$ws.Range("H136").Value = 4924.909
$ws.Range("I136").Value = 4691
$ws.Range("K136").Value = 14073
$ws.Range("M136").Value = -11523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42.433334
$ws.Range("I2").Value = 51.458332
$ws.Range("J2").Value = 6.3333335
$ws.Range("K2").Value = 308.749992
$ws.Range("L2").Value = 38.000001
$ws.Range("M2").Value = -195.749992
$ws.Range("N2").Value = -264.000001
$ws.Range("H92").Value = 916.5
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12962.25
$ws.Range("I70").Value = 4832.6665
$ws.Range("J70").Value = 17840
$ws.Range("K70").Value = 4832.6665
$ws.Range("L70").Value = 17840
$ws.Range("M70").Value = -4562.6665
$ws.Range("N70").Value = -18380
$ws.Range("H73").Value = 12962.25
$ws.Range("I73").Value = 4832.6665
$ws.Range("J73").Value = 17840
$ws.Range("K73").Value = 4832.6665
$ws.Range("L73").Value = 17840
$ws.Range("M73").Value = -3896.6665
$ws.Range("N73").Value = -19712
$ws.Range("H97").Value = 997.85187
$ws.Range("I97").Value = 766.0526
$ws.Range("K97").Value = 766.0526
$ws.Range("M97").Value = -270.0526
$ws.Range("H102").Value = 76924030
$ws.Range("J102").Value = 500000480
$ws.Range("L102").Value = 500000480
$ws.Range("N102").Value = -500003724
$ws.Range("H132").Value = 2766.3
$ws.Range("I132").Value = 2460.9333
$ws.Range("K132").Value = 7382.7999
$ws.Range("M132").Value = -4852.7999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4218
$ws.Range("I22").Value = 2984.6667
$ws.Range("K22").Value = 2984.6667
$ws.Range("M22").Value = -2689.6667
$ws.Range("H27").Value = 4218
$ws.Range("I27").Value = 2984.6667
$ws.Range("K27").Value = 2984.6667
$ws.Range("M27").Value = -2877.6667
$ws.Range("H30").Value = 400034.5
$ws.Range("I30").Value = 400034.5
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 400034.5
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -399926.5
$ws.Range("N30").ClearContents()
$ws.Range("H46").Value = 858.1667
$ws.Range("J46").Value = 916.6667
$ws.Range("L46").Value = 916.6667
$ws.Range("N46").Value = -1292.6667
$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2500
$ws.Range("N100").Value = -3582
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 5569.7144
$ws.Range("I122").Value = 4577.7
$ws.Range("K122").Value = 13733.1
$ws.Range("M122").Value = -11283.1
$ws.Range("H132").Value = 4358.5557
$ws.Range("I132").Value = 4227.7856
$ws.Range("K132").Value = 12683.3568
$ws.Range("M132").Value = -10153.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1832
$ws.Range("I96").Value = 750
$ws.Range("K96").Value = 750
$ws.Range("M96").Value = 623
$ws.Range("H118").Value = 45555
$ws.Range("J118").Value = 45555
$ws.Range("L118").Value = 45555
$ws.Range("N118").Value = -48869
